$d = $word.ActiveDocument

$replacements = @(
    @("232×6=", "756×8="),
    @("332×3=", "864×4="),
    @("873×3=", "328×7="),
    @("691×9=", "692×4="),
    @("792×2=", "117×3="),
    @("935×2=", "413×6="),
    @("719×4=", "614×5="),
    @("818×7=", "349×5="),
    @("335×6=", "187×5="),
    @("733×6=", "371×8="),
    @("229×6=", "908×5="),
    @("524×5=", "482×8="),
    @("565×9=", "288×9="),
    @("481×7=", "405×9="),
    @("141×8=", "860×7="),
    @("183×9=", "218×9="),
    @("603×7=", "704×9="),
    @("529×3=", "908×6="),
    @("751×8=", "910×9="),
    @("763×9=", "724×2="),
    @("220×5=", "147×6="),
    @("696×4=", "415×8="),
    @("578×5=", "458×2="),
    @("291×9=", "751×5="),
    @("281×8=", "370×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
